$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D9").Value = "이상한 AI BigData 대학원에서 고생 중이신 편입생도 받습니다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/transfer-students/#utm_source=rss&utm_medium=rss&utm_campaign=transfer-students"

$ws.Range("D50").Value = "데이터 기반 결정구조 탐색"
$ws.Range("E50").Value = "http://incredible.egloos.com/7515137"

$ws.Range("D51").Value = "[anaconda+python] spyder 5.0.0 버전 열자마자 나오는 에러 메시지(You have missing dependencies! # Mandatory: spyder_kernels)에 관하여"
$ws.Range("E51").Value = "https://bskyvision.com/1171"
